# Build site at 2022-09-26 16:07:08 UTC
# LOB1056.xlsx content refresh: the "Docentes responsaveis" value row (old
# row 13, label-less B/C pair holding the professor name) is removed from
# the sheet (shifting every row below it up by one), and several of the
# remaining label rows get their B/C values replaced with the text that the
# refreshed export now carries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stand-alone value row (old row 13: B/C = "8870322 - Fabiano
# Fernandes Bargos", no label in column A). Everything below shifts up by
# one row, which also takes care of the dimension/row-height changes.
$ws.Rows(13).Delete()

# --- Update the cells whose text content changed after the shift ---

# "Objetivos:" row keeps its label but now shows the docente text instead
# of the old Portuguese objectives paragraph.
$ws.Range("B10").Value = "8870322 - Fabiano Fernandes Bargos"
$ws.Range("C10").Value = "8870322 - Fabiano Fernandes Bargos"

# "Programa resumido:" row now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "Programa:" row now just shows the activation date.
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# "Metodo:" row now shows the docente text.
$ws.Range("B18").Value = "8870322 - Fabiano Fernandes Bargos"
$ws.Range("C18").Value = "8870322 - Fabiano Fernandes Bargos"

# "Criterio:" row now shows the evaluation-method paragraph.
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# "Norma de recuperação:" row now shows the "NF>=5,0." criterion text.
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# "Bibliografia:" row now shows the recovery-norm text (the long CHAPRA
# bibliography paragraph that used to be here is dropped entirely).
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."

# "Requisitos:" row (row 22) no longer carries a B/C value.
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
